$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($needle)) {
            return $i
        }
    }
    return -1
}

# Re-applies the "all zero" paragraph spacing/indent that the OM layer
# otherwise omits (as implied/default) once a paragraph has been
# rewritten via InsertXML.
function Fix-ZeroFormat($paraIndex) {
    $pp = $d.Paragraphs.Item($paraIndex)
    $pp.Format.SpaceBefore = 0
    $pp.Format.LeftIndent = 0
    $pp.Format.RightIndent = 0
    $pp.Format.FirstLineIndent = 0
}

# --- Common rPr (non-bold) block used throughout this vocab list ---
$plainRPr = '<w:rPr><w:rFonts w:ascii="맑은 고딕" w:hAnsi="맑은 고딕" w:cs="맑은 고딕" w:eastAsia="맑은 고딕"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="20"/><w:shd w:fill="auto" w:val="clear"/></w:rPr>'
$plainPPr = '<w:pPr><w:spacing w:before="0" w:after="200" w:line="276"/><w:ind w:right="0" w:left="0" w:firstLine="0"/><w:jc w:val="left"/>' + $plainRPr + '</w:pPr>'

# =====================================================================
# 1) "last time" / 마지막 -> merge the trailing translation run into the
#    run holding the English text + tabs (single run afterwards).
# =====================================================================
$idx = Find-ParagraphIndex("last time")
$p = $d.Paragraphs.Item($idx)
$xml = '<w:p ' + $wNs + '>' + $plainPPr + '<w:r>' + $plainRPr + `
    '<w:t xml:space="preserve">last time</w:t><w:tab/><w:tab/><w:tab/><w:tab/>' + `
    '<w:t xml:space="preserve">마지막</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)
Fix-ZeroFormat($idx)

# =====================================================================
# 2) "report" / 알리다. 신고하다. 전하다. -> same merge pattern.
# =====================================================================
$idx = Find-ParagraphIndex("report")
$p = $d.Paragraphs.Item($idx)
$xml = '<w:p ' + $wNs + '>' + $plainPPr + '<w:r>' + $plainRPr + `
    '<w:t xml:space="preserve">report</w:t><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/>' + `
    '<w:t xml:space="preserve">알리다. 신고하다. 전하다.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)
Fix-ZeroFormat($idx)

# =====================================================================
# 3) Insert 4 new vocabulary paragraphs right after "find out ... ".
# =====================================================================
$idx = Find-ParagraphIndex("find out")
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()

$newParasXml = `
    ('<w:p ' + $wNs + '>' + $plainPPr + `
        '<w:r>' + $plainRPr + '<w:t xml:space="preserve">organized </w:t><w:tab/><w:tab/><w:tab/><w:tab/></w:r>' + `
        '<w:r>' + $plainRPr + '<w:t xml:space="preserve">조직적인, 정리된, 체계적인</w:t></w:r>' + `
     '</w:p>') + `
    ('<w:p ' + $wNs + '>' + $plainPPr + `
        '<w:r>' + $plainRPr + `
        '<w:t xml:space="preserve">relevant </w:t><w:tab/><w:tab/><w:tab/><w:tab/>' + `
        '<w:t xml:space="preserve">관련있는, 적절한, 의의가 있는</w:t></w:r>' + `
     '</w:p>') + `
    ('<w:p ' + $wNs + '>' + $plainPPr + `
        '<w:r>' + $plainRPr + '<w:t xml:space="preserve">Moreover</w:t><w:tab/><w:tab/><w:tab/><w:tab/></w:r>' + `
        '<w:r>' + $plainRPr + '<w:t xml:space="preserve">게다가, 더욱이</w:t></w:r>' + `
     '</w:p>') + `
    ('<w:p ' + $wNs + '>' + $plainPPr + `
        '<w:r>' + $plainRPr + '<w:t xml:space="preserve">Impressive</w:t><w:tab/><w:tab/><w:tab/><w:tab/></w:r>' + `
        '<w:r>' + $plainRPr + '<w:t xml:space="preserve">인상적인, 인상 깊은</w:t></w:r>' + `
     '</w:p>')

$idx = Find-ParagraphIndex("find out")
$newP = $d.Paragraphs.Item($idx + 1)
$newP.Range.InsertXML($newParasXml)

for ($k = 0; $k -lt 4; $k++) {
    Fix-ZeroFormat($idx + 1 + $k)
}

# =====================================================================
# 4) "What's your card number?" / 카드번호가 뭡니까? -> merge the bold
#    translation run into the run holding the English text + tab.
# =====================================================================
$boldRPr = '<w:rPr><w:rFonts w:ascii="맑은 고딕" w:hAnsi="맑은 고딕" w:cs="맑은 고딕" w:eastAsia="맑은 고딕"/><w:b/><w:color w:val="auto"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="20"/><w:shd w:fill="auto" w:val="clear"/></w:rPr>'
$boldPPr = '<w:pPr><w:spacing w:before="0" w:after="200" w:line="276"/><w:ind w:right="0" w:left="0" w:firstLine="0"/><w:jc w:val="left"/>' + $boldRPr + '</w:pPr>'

$idx = Find-ParagraphIndex("What's your card number?")
$p = $d.Paragraphs.Item($idx)
$xml = '<w:p ' + $wNs + '>' + $boldPPr + '<w:r>' + $boldRPr + `
    '<w:t xml:space="preserve">What''s your card number?</w:t><w:tab/>' + `
    '<w:t xml:space="preserve">카드번호가 뭡니까?</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)
Fix-ZeroFormat($idx)
